$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.032.60"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.515.90"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.32"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.73"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "2.515.67"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").Value = "2.940.40"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.15"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "58.750.03"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "2.505.25"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.02"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.24"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.29"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("E24").Value = "  +5.52%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.58"
$ws.Range("E28").Value = "  -2.87%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "0.0₃0770"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.79"
$ws.Range("E32").Value = "  +4.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.17"
$ws.Range("E33").Value = "  +3.74%  "
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E35").Value = "  -4.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.53"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.79"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.21"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "279.49"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.604"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.88"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.98"
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.33"
$ws.Range("E51").Value = "  -1.13%  "
